$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44162
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 4667
$ws.Range("T2").Value = 1.5

$ws.Range("D3").Value = 44162
$ws.Range("M3").Value = 100
$ws.Range("O3").Value = 6500
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 4333

$ws.Range("D4").Value = 44516
$ws.Range("M4").Value = 80

$ws.Range("D5").Value = 44523
$ws.Range("M5").Value = 300

$ws.Range("D6").Value = 44169
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 5500
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 5750
$ws.Range("S6").Value = 3833

$ws.Range("D7").Value = 44176
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 5500
$ws.Range("S7").Value = 3667

$ws.Range("D8").Value = 44530
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 3600
$ws.Range("O8").Value = 3700
$ws.Range("P8").Value = 3650
$ws.Range("S8").Value = 3650

$ws.Range("D9").Value = 44159
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 6750
$ws.Range("S9").Value = 4500

$ws.Range("D10").Value = 44519
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 3700
$ws.Range("O10").Value = 3800
$ws.Range("P10").Value = 3750
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región del Maule"
$ws.Range("S10").Value = 3750
$ws.Range("T10").Value = 1
